$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exposed")
Write-Host $ws.Name
